$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values re-pulled from source data.
$ws.Range("F2").Value = -5
$ws.Range("F3").Value = -3
$ws.Range("F5").Value = -5
$ws.Range("F7").Value = -10
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = -3
$ws.Range("F19").Value = -2
$ws.Range("F21").Value = -3
$ws.Range("F23").Value = -2
$ws.Range("F31").Value = -5
$ws.Range("F35").Value = 2
$ws.Range("F37").Value = -1
